$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix swapped To (Country) / To (City) columns for existing flight rows
$ws.Cells.Item(2,2).Value = 'Aarhus'
$ws.Cells.Item(2,3).Value = 'Denmark'
$ws.Cells.Item(3,2).Value = 'Alghero'
$ws.Cells.Item(3,3).Value = 'Italy'
$ws.Cells.Item(4,2).Value = 'Almeria'
$ws.Cells.Item(4,3).Value = 'Spain'
$ws.Cells.Item(5,2).Value = 'Amsterdam'
$ws.Cells.Item(5,3).Value = 'Netherlands'
$ws.Cells.Item(6,2).Value = 'Amsterdam'
$ws.Cells.Item(6,3).Value = 'Netherlands'
$ws.Cells.Item(7,2).Value = 'Amsterdam'
$ws.Cells.Item(7,3).Value = 'Netherlands'
$ws.Cells.Item(8,2).Value = 'Antalya'
$ws.Cells.Item(8,3).Value = 'Turkey'
$ws.Cells.Item(9,2).Value = 'Athens'
$ws.Cells.Item(9,3).Value = 'Greece'
$ws.Cells.Item(10,2).Value = 'Baku'
$ws.Cells.Item(10,3).Value = 'Azerbaijan'
$ws.Cells.Item(11,2).Value = 'Barcelona'
$ws.Cells.Item(11,3).Value = 'Spain'
$ws.Cells.Item(12,2).Value = 'Barcelona'
$ws.Cells.Item(12,3).Value = 'Spain'
$ws.Cells.Item(13,2).Value = 'Barcelona'
$ws.Cells.Item(13,3).Value = 'Spain'
$ws.Cells.Item(14,2).Value = 'Bari'
$ws.Cells.Item(14,3).Value = 'Italy'
$ws.Cells.Item(15,2).Value = 'Basel/Mulhouse'
$ws.Cells.Item(15,3).Value = 'Switzerland'
$ws.Cells.Item(16,2).Value = 'Beijing'
$ws.Cells.Item(16,3).Value = 'China'
$ws.Cells.Item(17,2).Value = 'Beirut'
$ws.Cells.Item(17,3).Value = 'Lebanon'
$ws.Cells.Item(18,2).Value = 'Belgrade'
$ws.Cells.Item(18,3).Value = 'Serbia'
$ws.Cells.Item(19,2).Value = 'Bergen'
$ws.Cells.Item(19,3).Value = 'Norway'
$ws.Cells.Item(20,2).Value = 'Berlin/Tegel'
$ws.Cells.Item(20,3).Value = 'Germany'
$ws.Cells.Item(21,2).Value = 'Bilbao'
$ws.Cells.Item(21,3).Value = 'Spain'
$ws.Cells.Item(22,2).Value = 'Birmingham'
$ws.Cells.Item(22,3).Value = 'Great Britain'
$ws.Cells.Item(23,2).Value = 'Bologna'
$ws.Cells.Item(23,3).Value = 'Italy'
$ws.Cells.Item(24,2).Value = 'Bordeaux'
$ws.Cells.Item(24,3).Value = 'France'
$ws.Cells.Item(25,2).Value = 'Bourgas'
$ws.Cells.Item(25,3).Value = 'Bulgaria'
$ws.Cells.Item(26,2).Value = 'Bratislava'
$ws.Cells.Item(26,3).Value = 'Slovakia'
$ws.Cells.Item(27,2).Value = 'Bristol'
$ws.Cells.Item(27,3).Value = 'Great Britain'
$ws.Cells.Item(28,2).Value = 'Brussels/Charleroi'
$ws.Cells.Item(28,3).Value = 'Belgium'
$ws.Cells.Item(29,2).Value = 'Brussels/Zaventem'
$ws.Cells.Item(29,3).Value = 'Belgium'
$ws.Cells.Item(30,2).Value = 'Brussels/Zaventem'
$ws.Cells.Item(30,3).Value = 'Belgium'
$ws.Cells.Item(31,2).Value = 'Bucharest'
$ws.Cells.Item(31,3).Value = 'Romania'
$ws.Cells.Item(32,2).Value = 'Bucharest'
$ws.Cells.Item(32,3).Value = 'Romania'
$ws.Cells.Item(33,2).Value = 'Budapest'
$ws.Cells.Item(33,3).Value = 'Hungary'
$ws.Cells.Item(34,2).Value = 'Cagliari'
$ws.Cells.Item(34,3).Value = 'Italy'
$ws.Cells.Item(35,2).Value = 'Catania'
$ws.Cells.Item(35,3).Value = 'Italy'
$ws.Cells.Item(36,2).Value = 'Chania'
$ws.Cells.Item(36,3).Value = 'Greece'
$ws.Cells.Item(37,2).Value = 'Chengdu'
$ws.Cells.Item(37,3).Value = 'China'
$ws.Cells.Item(38,2).Value = 'Cologne/Bonn'
$ws.Cells.Item(38,3).Value = 'Germany'
$ws.Cells.Item(39,2).Value = 'Cologne/Bonn'
$ws.Cells.Item(39,3).Value = 'Germany'
$ws.Cells.Item(40,2).Value = 'Copenhagen'
$ws.Cells.Item(40,3).Value = 'Denmark'
$ws.Cells.Item(41,2).Value = 'Copenhagen'
$ws.Cells.Item(41,3).Value = 'Denmark'
$ws.Cells.Item(42,2).Value = 'Copenhagen'
$ws.Cells.Item(42,3).Value = 'Denmark'
$ws.Cells.Item(43,2).Value = 'Corfu'
$ws.Cells.Item(43,3).Value = 'Greece'
$ws.Cells.Item(44,2).Value = 'Dubai'
$ws.Cells.Item(44,3).Value = 'United Arab Emirates'
$ws.Cells.Item(45,2).Value = 'Dubai'
$ws.Cells.Item(45,3).Value = 'United Arab Emirates'
$ws.Cells.Item(46,2).Value = 'Dubai'
$ws.Cells.Item(46,3).Value = 'United Arab Emirates'
$ws.Cells.Item(47,2).Value = 'Dublin'
$ws.Cells.Item(47,3).Value = 'Ireland'
$ws.Cells.Item(48,2).Value = 'Dublin'
$ws.Cells.Item(48,3).Value = 'Ireland'
$ws.Cells.Item(49,2).Value = 'Dubrovnik'
$ws.Cells.Item(49,3).Value = 'Croatia'
$ws.Cells.Item(50,2).Value = 'Dusseldorf'
$ws.Cells.Item(50,3).Value = 'Germany'
$ws.Cells.Item(51,2).Value = 'Dusseldorf'
$ws.Cells.Item(51,3).Value = 'Germany'
$ws.Cells.Item(52,2).Value = 'Dusseldorf'
$ws.Cells.Item(52,3).Value = 'Germany'
$ws.Cells.Item(53,2).Value = 'East Midlands'
$ws.Cells.Item(53,3).Value = 'Great Britain'
$ws.Cells.Item(54,2).Value = 'Edinburgh'
$ws.Cells.Item(54,3).Value = 'Great Britain'
$ws.Cells.Item(55,2).Value = 'Eindhoven'
$ws.Cells.Item(55,3).Value = 'Netherlands'
$ws.Cells.Item(56,2).Value = 'Ekaterinburg'
$ws.Cells.Item(56,3).Value = 'Russia'
$ws.Cells.Item(57,2).Value = 'Ekaterinburg'
$ws.Cells.Item(57,3).Value = 'Russia'
$ws.Cells.Item(58,2).Value = 'Faro'
$ws.Cells.Item(58,3).Value = 'Portugal'
$ws.Cells.Item(59,2).Value = 'Frankfurt'
$ws.Cells.Item(59,3).Value = 'Germany'
$ws.Cells.Item(60,2).Value = 'Frankfurt'
$ws.Cells.Item(60,3).Value = 'Germany'
$ws.Cells.Item(61,2).Value = 'Fuerteventura'
$ws.Cells.Item(61,3).Value = 'Spain'
$ws.Cells.Item(62,2).Value = 'Geneva'
$ws.Cells.Item(62,3).Value = 'Switzerland'
$ws.Cells.Item(63,2).Value = 'Glasgow'
$ws.Cells.Item(63,3).Value = 'Great Britain'
$ws.Cells.Item(64,2).Value = 'Gothenburg'
$ws.Cells.Item(64,3).Value = 'Sweden'
$ws.Cells.Item(65,2).Value = 'Hamburg'
$ws.Cells.Item(65,3).Value = 'Germany'
$ws.Cells.Item(66,2).Value = 'Hamburg'
$ws.Cells.Item(66,3).Value = 'Germany'
$ws.Cells.Item(67,2).Value = 'Hamburg'
$ws.Cells.Item(67,3).Value = 'Germany'
$ws.Cells.Item(68,2).Value = 'Helsinki'
$ws.Cells.Item(68,3).Value = 'Finland'
$ws.Cells.Item(69,2).Value = 'Helsinki'
$ws.Cells.Item(69,3).Value = 'Finland'
$ws.Cells.Item(70,2).Value = 'Helsinki'
$ws.Cells.Item(70,3).Value = 'Finland'
$ws.Cells.Item(71,2).Value = 'Heraklion'
$ws.Cells.Item(71,3).Value = 'Greece'
$ws.Cells.Item(72,2).Value = 'Hévíz/Balaton'
$ws.Cells.Item(72,3).Value = 'Hungary'
$ws.Cells.Item(73,2).Value = 'Hurghada'
$ws.Cells.Item(73,3).Value = 'Egypt'
$ws.Cells.Item(74,2).Value = 'Hurghada'
$ws.Cells.Item(74,3).Value = 'Egypt'
$ws.Cells.Item(75,2).Value = 'Ibiza'
$ws.Cells.Item(75,3).Value = 'Spain'
$ws.Cells.Item(76,2).Value = 'Istanbul/Atatürk'
$ws.Cells.Item(76,3).Value = 'Turkey'
$ws.Cells.Item(77,2).Value = 'Istanbul/Sabiha Gokcen'
$ws.Cells.Item(77,3).Value = 'Turkey'
$ws.Cells.Item(78,2).Value = 'Kazan'
$ws.Cells.Item(78,3).Value = 'Russia'
$ws.Cells.Item(79,2).Value = 'Kefallinia'
$ws.Cells.Item(79,3).Value = 'Greece'
$ws.Cells.Item(80,2).Value = 'Keflavik'
$ws.Cells.Item(80,3).Value = 'Iceland'
$ws.Cells.Item(81,2).Value = 'Keflavik'
$ws.Cells.Item(81,3).Value = 'Iceland'
$ws.Cells.Item(82,2).Value = 'Kiev/Borispol'
$ws.Cells.Item(82,3).Value = 'Ukraine'
$ws.Cells.Item(83,2).Value = 'Kiev/Borispol'
$ws.Cells.Item(83,3).Value = 'Ukraine'
$ws.Cells.Item(84,2).Value = 'Kos'
$ws.Cells.Item(84,3).Value = 'Greece'
$ws.Cells.Item(85,2).Value = 'Kosice'
$ws.Cells.Item(85,3).Value = 'Slovakia'
$ws.Cells.Item(86,2).Value = 'Krasnodar'
$ws.Cells.Item(86,3).Value = 'Russia'
$ws.Cells.Item(87,2).Value = 'Lamezia Terme'
$ws.Cells.Item(87,3).Value = 'Italy'
$ws.Cells.Item(88,2).Value = 'Lanzarote'
$ws.Cells.Item(88,3).Value = 'Spain'
$ws.Cells.Item(89,2).Value = 'Larnaca'
$ws.Cells.Item(89,3).Value = 'Cyprus'
$ws.Cells.Item(90,2).Value = 'Las Palmas'
$ws.Cells.Item(90,3).Value = 'Spain'
$ws.Cells.Item(91,2).Value = 'Leeds'
$ws.Cells.Item(91,3).Value = 'Great Britain'
$ws.Cells.Item(92,2).Value = 'Linkoping'
$ws.Cells.Item(92,3).Value = 'Sweden'
$ws.Cells.Item(93,2).Value = 'Lisbon'
$ws.Cells.Item(93,3).Value = 'Portugal'
$ws.Cells.Item(94,2).Value = 'Lisbon'
$ws.Cells.Item(94,3).Value = 'Portugal'
$ws.Cells.Item(95,2).Value = 'Liverpool'
$ws.Cells.Item(95,3).Value = 'Great Britain'
$ws.Cells.Item(96,2).Value = 'Ljubljana'
$ws.Cells.Item(96,3).Value = 'Slovenia'
$ws.Cells.Item(97,2).Value = 'London/Gatwick'
$ws.Cells.Item(97,3).Value = 'Great Britain'
$ws.Cells.Item(98,2).Value = 'London/Gatwick'
$ws.Cells.Item(98,3).Value = 'Great Britain'
$ws.Cells.Item(99,2).Value = 'London/Heathrow'
$ws.Cells.Item(99,3).Value = 'Great Britain'
$ws.Cells.Item(100,2).Value = 'London/Luton'
$ws.Cells.Item(100,3).Value = 'Great Britain'
$ws.Cells.Item(101,2).Value = 'London/Southend'
$ws.Cells.Item(101,3).Value = 'Great Britain'
$ws.Cells.Item(102,2).Value = 'London/Stansted'
$ws.Cells.Item(102,3).Value = 'Great Britain'
$ws.Cells.Item(103,2).Value = 'London/Stansted'
$ws.Cells.Item(103,3).Value = 'Great Britain'
$ws.Cells.Item(104,2).Value = 'Luxembourg'
$ws.Cells.Item(104,3).Value = 'Luxembourg'
$ws.Cells.Item(105,2).Value = 'Lyon'
$ws.Cells.Item(105,3).Value = 'France'
$ws.Cells.Item(106,2).Value = 'Madeira/Funchal'
$ws.Cells.Item(106,3).Value = 'Portugal'
$ws.Cells.Item(107,2).Value = 'Madrid'
$ws.Cells.Item(107,3).Value = 'Spain'
$ws.Cells.Item(108,2).Value = 'Madrid'
$ws.Cells.Item(108,3).Value = 'Spain'
$ws.Cells.Item(109,2).Value = 'Malaga'
$ws.Cells.Item(109,3).Value = 'Spain'
$ws.Cells.Item(110,2).Value = 'Malta'
$ws.Cells.Item(110,3).Value = 'Malta'
$ws.Cells.Item(111,2).Value = 'Malta'
$ws.Cells.Item(111,3).Value = 'Malta'
$ws.Cells.Item(112,2).Value = 'Manchester'
$ws.Cells.Item(112,3).Value = 'Great Britain'
$ws.Cells.Item(113,2).Value = 'Manchester'
$ws.Cells.Item(113,3).Value = 'Great Britain'
$ws.Cells.Item(114,2).Value = 'Marsa Alam'
$ws.Cells.Item(114,3).Value = 'Egypt'
$ws.Cells.Item(115,2).Value = 'Marsa Alam'
$ws.Cells.Item(115,3).Value = 'Egypt'
$ws.Cells.Item(116,2).Value = 'Marseille'
$ws.Cells.Item(116,3).Value = 'France'
$ws.Cells.Item(117,2).Value = 'Menorca'
$ws.Cells.Item(117,3).Value = 'Spain'
$ws.Cells.Item(118,2).Value = 'Milan/Bergamo'
$ws.Cells.Item(118,3).Value = 'Italy'
$ws.Cells.Item(119,2).Value = 'Milan/Bergamo'
$ws.Cells.Item(119,3).Value = 'Italy'
$ws.Cells.Item(120,2).Value = 'Milan/Malpensa'
$ws.Cells.Item(120,3).Value = 'Italy'
$ws.Cells.Item(121,2).Value = 'Milan/Malpensa'
$ws.Cells.Item(121,3).Value = 'Italy'
$ws.Cells.Item(122,2).Value = 'Minsk'
$ws.Cells.Item(122,3).Value = 'Belarus'
$ws.Cells.Item(123,2).Value = 'Montreal'
$ws.Cells.Item(123,3).Value = 'Canada'
$ws.Cells.Item(124,2).Value = 'Moscow/Sheremetyevo'
$ws.Cells.Item(124,3).Value = 'Russia'
$ws.Cells.Item(125,2).Value = 'Moscow/Sheremetyevo'
$ws.Cells.Item(125,3).Value = 'Russia'
$ws.Cells.Item(126,2).Value = 'Moscow/Sheremetyevo'
$ws.Cells.Item(126,3).Value = 'Russia'
$ws.Cells.Item(127,2).Value = 'Munich'
$ws.Cells.Item(127,3).Value = 'Germany'
$ws.Cells.Item(128,2).Value = 'Nantes'
$ws.Cells.Item(128,3).Value = 'France'
$ws.Cells.Item(129,2).Value = 'Naples'
$ws.Cells.Item(129,3).Value = 'Italy'
$ws.Cells.Item(130,2).Value = 'Naples'
$ws.Cells.Item(130,3).Value = 'Italy'
$ws.Cells.Item(131,2).Value = 'Naples'
$ws.Cells.Item(131,3).Value = 'Italy'
$ws.Cells.Item(132,2).Value = 'New York/JFK'
$ws.Cells.Item(132,3).Value = 'USA'
$ws.Cells.Item(133,2).Value = 'Newcastle'
$ws.Cells.Item(133,3).Value = 'Great Britain'
$ws.Cells.Item(134,2).Value = 'Nice'
$ws.Cells.Item(134,3).Value = 'France'
$ws.Cells.Item(135,2).Value = 'Novosibirsk'
$ws.Cells.Item(135,3).Value = 'Russia'
$ws.Cells.Item(136,2).Value = 'Odessa'
$ws.Cells.Item(136,3).Value = 'Ukraine'
$ws.Cells.Item(137,2).Value = 'Olbia'
$ws.Cells.Item(137,3).Value = 'Italy'
$ws.Cells.Item(138,2).Value = 'Oslo'
$ws.Cells.Item(138,3).Value = 'Norway'
$ws.Cells.Item(139,2).Value = 'Ostrava'
$ws.Cells.Item(139,3).Value = 'Czech Republic'
$ws.Cells.Item(140,2).Value = 'Palma Mallorca'
$ws.Cells.Item(140,3).Value = 'Spain'
$ws.Cells.Item(141,2).Value = 'Paris/CDG'
$ws.Cells.Item(141,3).Value = 'France'
$ws.Cells.Item(142,2).Value = 'Paris/CDG'
$ws.Cells.Item(142,3).Value = 'France'
$ws.Cells.Item(143,2).Value = 'Paris/CDG'
$ws.Cells.Item(143,3).Value = 'France'
$ws.Cells.Item(144,2).Value = 'Paris/CDG'
$ws.Cells.Item(144,3).Value = 'France'
$ws.Cells.Item(145,2).Value = 'Paris/Orly'
$ws.Cells.Item(145,3).Value = 'France'
$ws.Cells.Item(146,2).Value = 'Pisa'
$ws.Cells.Item(146,3).Value = 'Italy'
$ws.Cells.Item(147,2).Value = 'Podgorica'
$ws.Cells.Item(147,3).Value = 'Monte Negro'
$ws.Cells.Item(148,2).Value = 'Porto'
$ws.Cells.Item(148,3).Value = 'Portugal'
$ws.Cells.Item(149,2).Value = 'Preveza'
$ws.Cells.Item(149,3).Value = 'Greece'
$ws.Cells.Item(150,2).Value = 'Radom'
$ws.Cells.Item(150,3).Value = 'Poland'
$ws.Cells.Item(151,2).Value = 'Ras Al Khaimah'
$ws.Cells.Item(151,3).Value = 'United Arab Emirates'
$ws.Cells.Item(152,2).Value = 'Rhodes'
$ws.Cells.Item(152,3).Value = 'Greece'
$ws.Cells.Item(153,2).Value = 'Riga'
$ws.Cells.Item(153,3).Value = 'Latvia'
$ws.Cells.Item(154,2).Value = 'Riyadh'
$ws.Cells.Item(154,3).Value = 'Saudi Arabia'
$ws.Cells.Item(155,2).Value = 'Rome/Ciampino'
$ws.Cells.Item(155,3).Value = 'Italy'
$ws.Cells.Item(156,2).Value = 'Rome/Fiumicino'
$ws.Cells.Item(156,3).Value = 'Italy'
$ws.Cells.Item(157,2).Value = 'Rome/Fiumicino'
$ws.Cells.Item(157,3).Value = 'Italy'
$ws.Cells.Item(158,2).Value = 'Rome/Fiumicino'
$ws.Cells.Item(158,3).Value = 'Italy'
$ws.Cells.Item(159,2).Value = 'Rostov on Don'
$ws.Cells.Item(159,3).Value = 'Russia'
$ws.Cells.Item(160,2).Value = 'Samara'
$ws.Cells.Item(160,3).Value = 'Russia'
$ws.Cells.Item(161,2).Value = 'Samos'
$ws.Cells.Item(161,3).Value = 'Greece'
$ws.Cells.Item(162,2).Value = 'Seoul/Incheon'
$ws.Cells.Item(162,3).Value = 'South Korea'
$ws.Cells.Item(163,2).Value = 'Seoul/Incheon'
$ws.Cells.Item(163,3).Value = 'South Korea'
$ws.Cells.Item(164,2).Value = 'Sevilla'
$ws.Cells.Item(164,3).Value = 'Spain'
$ws.Cells.Item(165,2).Value = 'Shanghai'
$ws.Cells.Item(165,3).Value = 'China'
$ws.Cells.Item(166,2).Value = 'Sharm El Sheikh'
$ws.Cells.Item(166,3).Value = 'Egypt'
$ws.Cells.Item(167,2).Value = 'Skopje'
$ws.Cells.Item(167,3).Value = 'Macedonia'
$ws.Cells.Item(168,2).Value = 'Sofia'
$ws.Cells.Item(168,3).Value = 'Bulgaria'
$ws.Cells.Item(169,2).Value = 'Split'
$ws.Cells.Item(169,3).Value = 'Croatia'
$ws.Cells.Item(170,2).Value = 'St Petersburg'
$ws.Cells.Item(170,3).Value = 'Russia'
$ws.Cells.Item(171,2).Value = 'St Petersburg'
$ws.Cells.Item(171,3).Value = 'Russia'
$ws.Cells.Item(172,2).Value = 'Stavanger'
$ws.Cells.Item(172,3).Value = 'Norway'
$ws.Cells.Item(173,2).Value = 'Stockholm/Arlanda'
$ws.Cells.Item(173,3).Value = 'Sweden'
$ws.Cells.Item(174,2).Value = 'Stockholm/Arlanda'
$ws.Cells.Item(174,3).Value = 'Sweden'
$ws.Cells.Item(175,2).Value = 'Stockholm/Arlanda'
$ws.Cells.Item(175,3).Value = 'Sweden'
$ws.Cells.Item(176,2).Value = 'Strasbourg'
$ws.Cells.Item(176,3).Value = 'France'
$ws.Cells.Item(177,2).Value = 'Tbilisi'
$ws.Cells.Item(177,3).Value = 'Georgia'
$ws.Cells.Item(178,2).Value = 'Tel Aviv'
$ws.Cells.Item(178,3).Value = 'Israel'
$ws.Cells.Item(179,2).Value = 'Tel Aviv'
$ws.Cells.Item(179,3).Value = 'Israel'
$ws.Cells.Item(180,2).Value = 'Tel Aviv'
$ws.Cells.Item(180,3).Value = 'Israel'
$ws.Cells.Item(181,2).Value = 'Tenerife'
$ws.Cells.Item(181,3).Value = 'Spain'
$ws.Cells.Item(182,2).Value = 'Thessaloniki'
$ws.Cells.Item(182,3).Value = 'Greece'
$ws.Cells.Item(183,2).Value = 'Tirana'
$ws.Cells.Item(183,3).Value = 'Albania'
$ws.Cells.Item(184,2).Value = 'Toronto'
$ws.Cells.Item(184,3).Value = 'Canada'
$ws.Cells.Item(185,2).Value = 'Toronto'
$ws.Cells.Item(185,3).Value = 'Canada'
$ws.Cells.Item(186,2).Value = 'Toulouse'
$ws.Cells.Item(186,3).Value = 'France'
$ws.Cells.Item(187,2).Value = 'Trapani'
$ws.Cells.Item(187,3).Value = 'Italy'
$ws.Cells.Item(188,2).Value = 'Tunis'
$ws.Cells.Item(188,3).Value = 'Tunisia'
$ws.Cells.Item(189,2).Value = 'Ufa'
$ws.Cells.Item(189,3).Value = 'Russia'
$ws.Cells.Item(190,2).Value = 'Valencia'
$ws.Cells.Item(190,3).Value = 'Spain'
$ws.Cells.Item(191,2).Value = 'Varna'
$ws.Cells.Item(191,3).Value = 'Bulgaria'
$ws.Cells.Item(192,2).Value = 'Vaxjo'
$ws.Cells.Item(192,3).Value = 'Sweden'
$ws.Cells.Item(193,2).Value = 'Venice/Marco Polo'
$ws.Cells.Item(193,3).Value = 'Italy'
$ws.Cells.Item(194,2).Value = 'Venice/Marco Polo'
$ws.Cells.Item(194,3).Value = 'Italy'
$ws.Cells.Item(195,2).Value = 'Venice/Marco Polo'
$ws.Cells.Item(195,3).Value = 'Italy'
$ws.Cells.Item(196,2).Value = 'Venice/Treviso'
$ws.Cells.Item(196,3).Value = 'Italy'
$ws.Cells.Item(197,2).Value = 'Verona'
$ws.Cells.Item(197,3).Value = 'Italy'
$ws.Cells.Item(198,2).Value = 'Vienna'
$ws.Cells.Item(198,3).Value = 'Austria'
$ws.Cells.Item(199,2).Value = 'Warsaw'
$ws.Cells.Item(199,3).Value = 'Poland'
$ws.Cells.Item(200,2).Value = 'Warsaw'
$ws.Cells.Item(200,3).Value = 'Poland'
$ws.Cells.Item(201,2).Value = 'Zagreb'
$ws.Cells.Item(201,3).Value = 'Croatia'
$ws.Cells.Item(202,2).Value = 'Zagreb'
$ws.Cells.Item(202,3).Value = 'Croatia'
$ws.Cells.Item(203,2).Value = 'Zakinthos'
$ws.Cells.Item(203,3).Value = 'Greece'
$ws.Cells.Item(204,2).Value = 'Zurich'
$ws.Cells.Item(204,3).Value = 'Switzerland'
$ws.Cells.Item(205,2).Value = 'Zurich'
$ws.Cells.Item(205,3).Value = 'Switzerland'

# Append new flight rows for Brno and Ostrava
$ws.Cells.Item(206,1).Value = 'Brno'
$ws.Cells.Item(206,2).Value = 'ANTALYA'
$ws.Cells.Item(206,3).Value = ''
$ws.Cells.Item(206,4).Value = 'Travel Service'
$ws.Cells.Item(207,1).Value = 'Brno'
$ws.Cells.Item(207,2).Value = 'KOS'
$ws.Cells.Item(207,3).Value = ''
$ws.Cells.Item(207,4).Value = 'Travel Service'
$ws.Cells.Item(208,1).Value = 'Brno'
$ws.Cells.Item(208,2).Value = 'BURGAS'
$ws.Cells.Item(208,3).Value = ''
$ws.Cells.Item(208,4).Value = 'Travel Service'
$ws.Cells.Item(209,1).Value = 'Brno'
$ws.Cells.Item(209,2).Value = 'EINDHOVEN'
$ws.Cells.Item(209,3).Value = ''
$ws.Cells.Item(209,4).Value = 'Wizz Air'
$ws.Cells.Item(210,1).Value = 'Brno'
$ws.Cells.Item(210,2).Value = 'MUNICH'
$ws.Cells.Item(210,3).Value = ''
$ws.Cells.Item(210,4).Value = 'bmi regional'
$ws.Cells.Item(211,1).Value = 'Brno'
$ws.Cells.Item(211,2).Value = 'MUNICH'
$ws.Cells.Item(211,3).Value = ''
$ws.Cells.Item(211,4).Value = 'Lufthansa'
$ws.Cells.Item(212,1).Value = 'Brno'
$ws.Cells.Item(212,2).Value = 'ZAKYNTHOS'
$ws.Cells.Item(212,3).Value = ''
$ws.Cells.Item(212,4).Value = 'Travel Service'
$ws.Cells.Item(213,1).Value = 'Brno'
$ws.Cells.Item(213,2).Value = 'MONASTIR / OSTRAVA'
$ws.Cells.Item(213,3).Value = ''
$ws.Cells.Item(213,4).Value = 'Tunisair'
$ws.Cells.Item(214,1).Value = 'Brno'
$ws.Cells.Item(214,2).Value = 'LONDON STANSTED'
$ws.Cells.Item(214,3).Value = ''
$ws.Cells.Item(214,4).Value = 'Ryanair'
$ws.Cells.Item(215,1).Value = 'Brno'
$ws.Cells.Item(215,2).Value = 'LAMEZIA TERME'
$ws.Cells.Item(215,3).Value = ''
$ws.Cells.Item(215,4).Value = 'Travel Service'
$ws.Cells.Item(216,1).Value = 'Brno'
$ws.Cells.Item(216,2).Value = 'VARNA'
$ws.Cells.Item(216,3).Value = ''
$ws.Cells.Item(216,4).Value = 'Travel Service'
$ws.Cells.Item(217,1).Value = 'Brno'
$ws.Cells.Item(217,2).Value = 'RHODOS'
$ws.Cells.Item(217,3).Value = ''
$ws.Cells.Item(217,4).Value = 'Travel Service'
$ws.Cells.Item(218,1).Value = 'Brno'
$ws.Cells.Item(218,2).Value = 'LONDON LUTON'
$ws.Cells.Item(218,3).Value = ''
$ws.Cells.Item(218,4).Value = 'Wizz Air'
$ws.Cells.Item(219,1).Value = 'Brno'
$ws.Cells.Item(219,2).Value = 'HERAKLION'
$ws.Cells.Item(219,3).Value = ''
$ws.Cells.Item(219,4).Value = 'Travel Service'
$ws.Cells.Item(220,1).Value = 'Brno'
$ws.Cells.Item(220,2).Value = 'PALMA DE MALLORCA'
$ws.Cells.Item(220,3).Value = ''
$ws.Cells.Item(220,4).Value = 'Travel Service'
$ws.Cells.Item(221,1).Value = 'Brno'
$ws.Cells.Item(221,2).Value = 'ERCAN / ANTALYA'
$ws.Cells.Item(221,3).Value = ''
$ws.Cells.Item(221,4).Value = 'Tailwind Airlines'
$ws.Cells.Item(222,1).Value = 'Brno'
$ws.Cells.Item(222,2).Value = 'CORFU'
$ws.Cells.Item(222,3).Value = ''
$ws.Cells.Item(222,4).Value = 'Travel Service'
$ws.Cells.Item(223,1).Value = 'Brno'
$ws.Cells.Item(223,2).Value = 'ENFIDHA'
$ws.Cells.Item(223,3).Value = ''
$ws.Cells.Item(223,4).Value = 'Tunisair'
$ws.Cells.Item(224,1).Value = 'Brno'
$ws.Cells.Item(224,2).Value = 'DJERBA / OSTRAVA'
$ws.Cells.Item(224,3).Value = ''
$ws.Cells.Item(224,4).Value = 'Tunisair'
$ws.Cells.Item(225,1).Value = 'Brno'
$ws.Cells.Item(225,2).Value = 'DJERBA / ENFIDHA'
$ws.Cells.Item(225,3).Value = ''
$ws.Cells.Item(225,4).Value = 'Travel Service'
$ws.Cells.Item(226,1).Value = 'Brno'
$ws.Cells.Item(226,2).Value = 'HURGHADA'
$ws.Cells.Item(226,3).Value = ''
$ws.Cells.Item(226,4).Value = 'Travel Service'
$ws.Cells.Item(227,1).Value = 'Brno'
$ws.Cells.Item(227,2).Value = 'PODGORICA'
$ws.Cells.Item(227,3).Value = ''
$ws.Cells.Item(227,4).Value = 'Travel Service'
$ws.Cells.Item(228,1).Value = 'Brno'
$ws.Cells.Item(228,2).Value = 'MARSA ALAM'
$ws.Cells.Item(228,3).Value = ''
$ws.Cells.Item(228,4).Value = 'Travel Service'
$ws.Cells.Item(229,1).Value = 'Brno'
$ws.Cells.Item(229,2).Value = 'VARNA / OSTRAVA'
$ws.Cells.Item(229,3).Value = ''
$ws.Cells.Item(229,4).Value = 'Travel Service'
$ws.Cells.Item(230,1).Value = 'Brno'
$ws.Cells.Item(230,2).Value = 'TENERIFE'
$ws.Cells.Item(230,3).Value = ''
$ws.Cells.Item(230,4).Value = 'Travel Service'
$ws.Cells.Item(231,1).Value = 'Brno'
$ws.Cells.Item(231,2).Value = 'THESSALONIKI'
$ws.Cells.Item(231,3).Value = ''
$ws.Cells.Item(231,4).Value = 'Travel Service'
$ws.Cells.Item(232,1).Value = 'Brno'
$ws.Cells.Item(232,2).Value = 'KAVALA / OSTRAVA'
$ws.Cells.Item(232,3).Value = ''
$ws.Cells.Item(232,4).Value = 'Travel Service'
$ws.Cells.Item(233,1).Value = 'Brno'
$ws.Cells.Item(233,2).Value = 'ALMERIA / OSTRAVA'
$ws.Cells.Item(233,3).Value = ''
$ws.Cells.Item(233,4).Value = 'Travel Service'
$ws.Cells.Item(234,1).Value = 'Brno'
$ws.Cells.Item(234,2).Value = 'PREVEZA'
$ws.Cells.Item(234,3).Value = ''
$ws.Cells.Item(234,4).Value = 'Travel Service'
$ws.Cells.Item(235,1).Value = 'Ostrava'
$ws.Cells.Item(235,2).Value = 'Antalya, Antayla (AYT)'
$ws.Cells.Item(235,3).Value = ''
$ws.Cells.Item(235,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(236,1).Value = 'Ostrava'
$ws.Cells.Item(236,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(236,3).Value = ''
$ws.Cells.Item(236,4).Value = 'CZECH AIRLINES (CSA)'
$ws.Cells.Item(237,1).Value = 'Ostrava'
$ws.Cells.Item(237,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(237,3).Value = ''
$ws.Cells.Item(237,4).Value = 'KLM ROYAL DUTCH AIRLINES'
$ws.Cells.Item(238,1).Value = 'Ostrava'
$ws.Cells.Item(238,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(238,3).Value = ''
$ws.Cells.Item(238,4).Value = 'DELTA AIR LINES'
$ws.Cells.Item(239,1).Value = 'Ostrava'
$ws.Cells.Item(239,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(239,3).Value = ''
$ws.Cells.Item(239,4).Value = 'KOREAN AIR'
$ws.Cells.Item(240,1).Value = 'Ostrava'
$ws.Cells.Item(240,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(240,3).Value = ''
$ws.Cells.Item(240,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(241,1).Value = 'Ostrava'
$ws.Cells.Item(241,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(241,3).Value = ''
$ws.Cells.Item(241,4).Value = 'TAROM ROMANIAN AIRLINES'
$ws.Cells.Item(242,1).Value = 'Ostrava'
$ws.Cells.Item(242,2).Value = 'Burgas, Burgas Airport (BOJ)'
$ws.Cells.Item(242,3).Value = ''
$ws.Cells.Item(242,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(243,1).Value = 'Ostrava'
$ws.Cells.Item(243,2).Value = 'Kos, Kos Island International Airport (KGS)'
$ws.Cells.Item(243,3).Value = ''
$ws.Cells.Item(243,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(244,1).Value = 'Ostrava'
$ws.Cells.Item(244,2).Value = 'Varna, Varna Airport (VAR)'
$ws.Cells.Item(244,3).Value = ''
$ws.Cells.Item(244,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(245,1).Value = 'Ostrava'
$ws.Cells.Item(245,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(245,3).Value = ''
$ws.Cells.Item(245,4).Value = 'CZECH AIRLINES (CSA)'
$ws.Cells.Item(246,1).Value = 'Ostrava'
$ws.Cells.Item(246,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(246,3).Value = ''
$ws.Cells.Item(246,4).Value = 'KLM ROYAL DUTCH AIRLINES'
$ws.Cells.Item(247,1).Value = 'Ostrava'
$ws.Cells.Item(247,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(247,3).Value = ''
$ws.Cells.Item(247,4).Value = 'KOREAN AIR'
$ws.Cells.Item(248,1).Value = 'Ostrava'
$ws.Cells.Item(248,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(248,3).Value = ''
$ws.Cells.Item(248,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(249,1).Value = 'Ostrava'
$ws.Cells.Item(249,2).Value = 'London, Stansted (STN)'
$ws.Cells.Item(249,3).Value = ''
$ws.Cells.Item(249,4).Value = 'RYANAIR'
$ws.Cells.Item(250,1).Value = 'Ostrava'
$ws.Cells.Item(250,2).Value = 'MONASTIR'
$ws.Cells.Item(250,3).Value = ''
$ws.Cells.Item(251,1).Value = 'Ostrava'
$ws.Cells.Item(251,2).Value = 'Rhodes, Diagoras Airport (RHO)'
$ws.Cells.Item(251,3).Value = ''
$ws.Cells.Item(251,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(252,1).Value = 'Ostrava'
$ws.Cells.Item(252,2).Value = 'Crete / Heraklion, N. Kazantzakis Apt. (HER)'
$ws.Cells.Item(252,3).Value = ''
$ws.Cells.Item(252,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(253,1).Value = 'Ostrava'
$ws.Cells.Item(253,2).Value = 'Burgas, Burgas Airport (BOJ)'
$ws.Cells.Item(253,3).Value = ''
$ws.Cells.Item(253,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(254,1).Value = 'Ostrava'
$ws.Cells.Item(254,2).Value = 'Mallorca, Palma de Mallorca (PMI)'
$ws.Cells.Item(254,3).Value = ''
$ws.Cells.Item(254,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(255,1).Value = 'Ostrava'
$ws.Cells.Item(255,2).Value = 'Rhodes, Diagoras Airport (RHO)'
$ws.Cells.Item(255,3).Value = ''
$ws.Cells.Item(255,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(256,1).Value = 'Ostrava'
$ws.Cells.Item(256,2).Value = 'Milan / Bergamo, Milan Bergamo Airport'
$ws.Cells.Item(256,3).Value = ''
$ws.Cells.Item(256,4).Value = 'RYANAIR'
$ws.Cells.Item(257,1).Value = 'Ostrava'
$ws.Cells.Item(257,2).Value = 'Corfu / Kerkyra, I. Kapodistrias (CFU)'
$ws.Cells.Item(257,3).Value = ''
$ws.Cells.Item(257,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(258,1).Value = 'Ostrava'
$ws.Cells.Item(258,2).Value = 'Djerba'
$ws.Cells.Item(258,3).Value = ''
$ws.Cells.Item(259,1).Value = 'Ostrava'
$ws.Cells.Item(259,2).Value = 'Zakynthos, Zakinthos Is (ZTH)'
$ws.Cells.Item(259,3).Value = ''
$ws.Cells.Item(259,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(260,1).Value = 'Ostrava'
$ws.Cells.Item(260,2).Value = 'Split'
$ws.Cells.Item(260,3).Value = ''
$ws.Cells.Item(260,4).Value = 'CZECH AIRLINES (CSA)'
$ws.Cells.Item(261,1).Value = 'Ostrava'
$ws.Cells.Item(261,2).Value = 'Hurghada, Hurghada (HRG)'
$ws.Cells.Item(261,3).Value = ''
$ws.Cells.Item(261,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(262,1).Value = 'Ostrava'
$ws.Cells.Item(262,2).Value = 'Burgas, Burgas Airport (BOJ)'
$ws.Cells.Item(262,3).Value = ''
$ws.Cells.Item(262,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(263,1).Value = 'Ostrava'
$ws.Cells.Item(263,2).Value = 'Podgorica'
$ws.Cells.Item(263,3).Value = ''
$ws.Cells.Item(263,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(264,1).Value = 'Ostrava'
$ws.Cells.Item(264,2).Value = 'Burgas, Burgas Airport (BOJ)'
$ws.Cells.Item(264,3).Value = ''
$ws.Cells.Item(264,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(265,1).Value = 'Ostrava'
$ws.Cells.Item(265,2).Value = 'Rotterdam'
$ws.Cells.Item(265,3).Value = ''
$ws.Cells.Item(266,1).Value = 'Ostrava'
$ws.Cells.Item(266,2).Value = 'Marsa Alam, Marsa Alam (RMF)'
$ws.Cells.Item(266,3).Value = ''
$ws.Cells.Item(266,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(267,1).Value = 'Ostrava'
$ws.Cells.Item(267,2).Value = 'Varna, Varna Airport (VAR)'
$ws.Cells.Item(267,3).Value = ''
$ws.Cells.Item(267,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(268,1).Value = 'Ostrava'
$ws.Cells.Item(268,2).Value = 'Warsawa'
$ws.Cells.Item(268,3).Value = ''
$ws.Cells.Item(268,4).Value = 'CZECH AIRLINES (CSA)'
$ws.Cells.Item(269,1).Value = 'Ostrava'
$ws.Cells.Item(269,2).Value = 'Warsawa'
$ws.Cells.Item(269,3).Value = ''
$ws.Cells.Item(269,4).Value = 'KLM ROYAL DUTCH AIRLINES'
$ws.Cells.Item(270,1).Value = 'Ostrava'
$ws.Cells.Item(270,2).Value = 'Warsawa'
$ws.Cells.Item(270,3).Value = ''
$ws.Cells.Item(270,4).Value = 'DELTA AIR LINES'
$ws.Cells.Item(271,1).Value = 'Ostrava'
$ws.Cells.Item(271,2).Value = 'Warsawa'
$ws.Cells.Item(271,3).Value = ''
$ws.Cells.Item(271,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(272,1).Value = 'Ostrava'
$ws.Cells.Item(272,2).Value = 'Warsawa'
$ws.Cells.Item(272,3).Value = ''
$ws.Cells.Item(272,4).Value = 'TAROM ROMANIAN AIRLINES'
$ws.Cells.Item(273,1).Value = 'Ostrava'
$ws.Cells.Item(273,2).Value = 'KAVALA'
$ws.Cells.Item(273,3).Value = ''
$ws.Cells.Item(273,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(274,1).Value = 'Ostrava'
$ws.Cells.Item(274,2).Value = 'ALMERIA'
$ws.Cells.Item(274,3).Value = ''
$ws.Cells.Item(274,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(275,1).Value = 'Ostrava'
$ws.Cells.Item(275,2).Value = 'Prague, Václav Havel Airport Prague (PRG)'
$ws.Cells.Item(275,3).Value = ''
$ws.Cells.Item(275,4).Value = 'TRAVEL SERVICE / SMARTWINGS'
$ws.Cells.Item(276,1).Value = 'Ostrava'
$ws.Cells.Item(276,2).Value = 'Podgorica'
$ws.Cells.Item(276,3).Value = ''
$ws.Cells.Item(276,4).Value = 'TRAVEL SERVICE / SMARTWINGS'

# Update selection to match saved view state
$ws.Range("A2").Select()
